$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two data values that changed on row 2
$ws.Range("A2").Value = 435435455
$ws.Range("C2").Value = 997650

# Bump the duplicate-values conditional formatting rule's priority (1 -> 2)
$range = $ws.Range("E2")
$fc = $range.FormatConditions.Item(1)
$fc.Priority = 2

# Move the active cell / selection to C11
$ws.Range("C11").Select()

$wb.Save()
